$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# This edit re-sorts the three status sheets (Overview, zh-cn, de-de) so the
# row describing file "198542c6-4d12-4bf5-824b-0062f0c832fd" - which just
# became "Ready for handoff" - moves from the top data row down to the
# bottom, while the two "Handed back" rows (ffffc5d5a3c0..., ffffff791f7ac9...)
# shift up to take its place. Along with the move, the 198542c6 row gets a
# refreshed handoff status + handoff timestamp on every sheet.
# ---------------------------------------------------------------------------

# ========================= Sheet "Overview" =========================
$wsO = $wb.Worksheets.Item("Overview")

$wsO.Range("A2").Value = "ffffc5d5a3c0-02f9-4907-ab3a-fb673b2c366d.md"
$wsO.Range("B2").Value = "Handed back: in sync with en-US"
$wsO.Range("C2").Value = "Handed back: in sync with en-US"
$wsO.Range("D2").Value = "2016-03-22 17:12:59"

$wsO.Range("A3").Value = "ffffff791f7ac9-e4ac-4cdf-acf3-2ad41c40a99e.md"
$wsO.Range("B3").Value = "Handed back: in sync with en-US"
$wsO.Range("C3").Value = "Handed back: in sync with en-US"
$wsO.Range("D3").Value = "2016-03-22 17:12:59"

$wsO.Range("A4").Value = "198542c6-4d12-4bf5-824b-0062f0c832fd.md"
$wsO.Range("B4").Value = "Ready for handoff"
$wsO.Range("C4").Value = "Ready for handoff"
$wsO.Range("D4").Value = "2016-03-22 17:17:23"

$wsO.Hyperlinks.Delete()
$wsO.Hyperlinks.Add($wsO.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/85732830159f30e0569f744e3032dbff5009f940/e2e/ffffc5d5a3c0-02f9-4907-ab3a-fb673b2c366d.md", "", "", "ffffc5d5a3c0-02f9-4907-ab3a-fb673b2c366d.md") | Out-Null
$wsO.Hyperlinks.Add($wsO.Range("A3"), "https://github.com/OpenLocalizationTest/oltest/blob/85732830159f30e0569f744e3032dbff5009f940/e2e/ffffff791f7ac9-e4ac-4cdf-acf3-2ad41c40a99e.md", "", "", "ffffff791f7ac9-e4ac-4cdf-acf3-2ad41c40a99e.md") | Out-Null
$wsO.Hyperlinks.Add($wsO.Range("A4"), "https://github.com/OpenLocalizationTest/oltest/blob/85732830159f30e0569f744e3032dbff5009f940/e2e/198542c6-4d12-4bf5-824b-0062f0c832fd.md", "", "", "198542c6-4d12-4bf5-824b-0062f0c832fd.md") | Out-Null

# ========================= Sheet "zh-cn" =========================
$wsZ = $wb.Worksheets.Item("zh-cn")

$wsZ.Range("A2").Value = "ffffc5d5a3c0-02f9-4907-ab3a-fb673b2c366d.md"
$wsZ.Range("B2").Value = ".md"
$wsZ.Range("C2").Value = "Handed back: in sync with en-US"
$wsZ.Range("D2").Value = "59db17ff-0d89-493f-a4fb-64bdf414a197.bac44bfe36bffe9cc476143af072f2ebdc47248a.zh-cn.xlf"
$wsZ.Range("E2").Value = "2016-03-22 17:12:55"
$wsZ.Range("F2").Value = "59db17ff-0d89-493f-a4fb-64bdf414a197.md"
$wsZ.Range("G2").Value = "59db17ff-0d89-493f-a4fb-64bdf414a197.bac44bfe36bffe9cc476143af072f2ebdc47248a.zh-cn.xlf"
$wsZ.Range("H2").Value = "2016-03-22 17:13:23"
$wsZ.Range("J2").Value = "Include"

$wsZ.Range("A3").Value = "ffffff791f7ac9-e4ac-4cdf-acf3-2ad41c40a99e.md"
$wsZ.Range("B3").Value = ".md"
$wsZ.Range("C3").Value = "Handed back: in sync with en-US"
$wsZ.Range("D3").Value = "59db17ff-0d89-493f-a4fb-64bdf414a197.bac44bfe36bffe9cc476143af072f2ebdc47248a.zh-cn.xlf"
$wsZ.Range("E3").Value = "2016-03-22 17:12:55"
$wsZ.Range("F3").Value = "59db17ff-0d89-493f-a4fb-64bdf414a197.md"
$wsZ.Range("G3").Value = "59db17ff-0d89-493f-a4fb-64bdf414a197.bac44bfe36bffe9cc476143af072f2ebdc47248a.zh-cn.xlf"
$wsZ.Range("H3").Value = "2016-03-22 17:13:23"
$wsZ.Range("J3").Value = "Include"

$wsZ.Range("A4").Value = "198542c6-4d12-4bf5-824b-0062f0c832fd.md"
$wsZ.Range("B4").Value = ".md"
$wsZ.Range("C4").Value = "Ready for handoff"
$wsZ.Range("D4").Value = "198542c6-4d12-4bf5-824b-0062f0c832fd.fa0c9f31959e458e7e6c4b18a2ee6277848d5c38.zh-cn.xlf"
$wsZ.Range("E4").Value = "2016-03-22 17:17:18"
$wsZ.Range("F4").Value = "198542c6-4d12-4bf5-824b-0062f0c832fd.md"
$wsZ.Range("G4").Value = "198542c6-4d12-4bf5-824b-0062f0c832fd.fa0c9f31959e458e7e6c4b18a2ee6277848d5c38.zh-cn.xlf"
$wsZ.Range("H4").Value = "2016-03-22 17:16:21"
$wsZ.Range("J4").Value = "Include"

$wsZ.Hyperlinks.Delete()
$wsZ.Hyperlinks.Add($wsZ.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/85732830159f30e0569f744e3032dbff5009f940/e2e/ffffc5d5a3c0-02f9-4907-ab3a-fb673b2c366d.md", "", "", "ffffc5d5a3c0-02f9-4907-ab3a-fb673b2c366d.md") | Out-Null
$wsZ.Hyperlinks.Add($wsZ.Range("D2"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/dfc61ff083b488bac160d5f280d22a652d52b7b0/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/59db17ff-0d89-493f-a4fb-64bdf414a197.bac44bfe36bffe9cc476143af072f2ebdc47248a.zh-cn.xlf", "", "", "59db17ff-0d89-493f-a4fb-64bdf414a197.bac44bfe36bffe9cc476143af072f2ebdc47248a.zh-cn.xlf") | Out-Null
$wsZ.Hyperlinks.Add($wsZ.Range("F2"), "https://github.com/OpenLocalizationTestOrg/oltest.zh-cn/blob/18ca270076bb2e209fa74fa1e19e73d876ad8138/e2e/59db17ff-0d89-493f-a4fb-64bdf414a197.md", "", "", "59db17ff-0d89-493f-a4fb-64bdf414a197.md") | Out-Null
$wsZ.Hyperlinks.Add($wsZ.Range("G2"), "https://github.com/OpenLocalizationTestOrg/olhandback/blob/df3d752db5bcafa12b5e329bbb94b79b5d85b08c/ol-handback/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/59db17ff-0d89-493f-a4fb-64bdf414a197.bac44bfe36bffe9cc476143af072f2ebdc47248a.zh-cn.xlf", "", "", "59db17ff-0d89-493f-a4fb-64bdf414a197.bac44bfe36bffe9cc476143af072f2ebdc47248a.zh-cn.xlf") | Out-Null

$wsZ.Hyperlinks.Add($wsZ.Range("A3"), "https://github.com/OpenLocalizationTest/oltest/blob/85732830159f30e0569f744e3032dbff5009f940/e2e/ffffff791f7ac9-e4ac-4cdf-acf3-2ad41c40a99e.md", "", "", "ffffff791f7ac9-e4ac-4cdf-acf3-2ad41c40a99e.md") | Out-Null
$wsZ.Hyperlinks.Add($wsZ.Range("D3"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/dfc61ff083b488bac160d5f280d22a652d52b7b0/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/59db17ff-0d89-493f-a4fb-64bdf414a197.bac44bfe36bffe9cc476143af072f2ebdc47248a.zh-cn.xlf", "", "", "59db17ff-0d89-493f-a4fb-64bdf414a197.bac44bfe36bffe9cc476143af072f2ebdc47248a.zh-cn.xlf") | Out-Null
$wsZ.Hyperlinks.Add($wsZ.Range("F3"), "https://github.com/OpenLocalizationTestOrg/oltest.zh-cn/blob/18ca270076bb2e209fa74fa1e19e73d876ad8138/e2e/59db17ff-0d89-493f-a4fb-64bdf414a197.md", "", "", "59db17ff-0d89-493f-a4fb-64bdf414a197.md") | Out-Null
$wsZ.Hyperlinks.Add($wsZ.Range("G3"), "https://github.com/OpenLocalizationTestOrg/olhandback/blob/df3d752db5bcafa12b5e329bbb94b79b5d85b08c/ol-handback/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/59db17ff-0d89-493f-a4fb-64bdf414a197.bac44bfe36bffe9cc476143af072f2ebdc47248a.zh-cn.xlf", "", "", "59db17ff-0d89-493f-a4fb-64bdf414a197.bac44bfe36bffe9cc476143af072f2ebdc47248a.zh-cn.xlf") | Out-Null

$wsZ.Hyperlinks.Add($wsZ.Range("A4"), "https://github.com/OpenLocalizationTest/oltest/blob/85732830159f30e0569f744e3032dbff5009f940/e2e/198542c6-4d12-4bf5-824b-0062f0c832fd.md", "", "", "198542c6-4d12-4bf5-824b-0062f0c832fd.md") | Out-Null
$wsZ.Hyperlinks.Add($wsZ.Range("D4"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/d1417235770bc0d2a6cef0e7c5e871df08a7502d/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/198542c6-4d12-4bf5-824b-0062f0c832fd.fa0c9f31959e458e7e6c4b18a2ee6277848d5c38.zh-cn.xlf", "", "", "198542c6-4d12-4bf5-824b-0062f0c832fd.fa0c9f31959e458e7e6c4b18a2ee6277848d5c38.zh-cn.xlf") | Out-Null
$wsZ.Hyperlinks.Add($wsZ.Range("F4"), "https://github.com/OpenLocalizationTestOrg/oltest.zh-cn/blob/82e26a8aa39123cb26ef0c18bf5380753a8ef30b/e2e/198542c6-4d12-4bf5-824b-0062f0c832fd.md", "", "", "198542c6-4d12-4bf5-824b-0062f0c832fd.md") | Out-Null
$wsZ.Hyperlinks.Add($wsZ.Range("G4"), "https://github.com/OpenLocalizationTestOrg/olhandback/blob/55de8e549bdea43e6a82e19b6856eded8ac07d1a/ol-handback/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/198542c6-4d12-4bf5-824b-0062f0c832fd.fa0c9f31959e458e7e6c4b18a2ee6277848d5c38.zh-cn.xlf", "", "", "198542c6-4d12-4bf5-824b-0062f0c832fd.fa0c9f31959e458e7e6c4b18a2ee6277848d5c38.zh-cn.xlf") | Out-Null

# ========================= Sheet "de-de" =========================
$wsD = $wb.Worksheets.Item("de-de")

$wsD.Range("A2").Value = "ffffc5d5a3c0-02f9-4907-ab3a-fb673b2c366d.md"
$wsD.Range("B2").Value = ".md"
$wsD.Range("C2").Value = "Handed back: in sync with en-US"
$wsD.Range("D2").Value = "59db17ff-0d89-493f-a4fb-64bdf414a197.bac44bfe36bffe9cc476143af072f2ebdc47248a.de-de.xlf"
$wsD.Range("E2").Value = "2016-03-22 17:12:59"
$wsD.Range("F2").Value = "59db17ff-0d89-493f-a4fb-64bdf414a197.md"
$wsD.Range("G2").Value = "59db17ff-0d89-493f-a4fb-64bdf414a197.bac44bfe36bffe9cc476143af072f2ebdc47248a.de-de.xlf"
$wsD.Range("H2").Value = "2016-03-22 17:13:30"
$wsD.Range("J2").Value = "Include"

$wsD.Range("A3").Value = "ffffff791f7ac9-e4ac-4cdf-acf3-2ad41c40a99e.md"
$wsD.Range("B3").Value = ".md"
$wsD.Range("C3").Value = "Handed back: in sync with en-US"
$wsD.Range("D3").Value = "59db17ff-0d89-493f-a4fb-64bdf414a197.bac44bfe36bffe9cc476143af072f2ebdc47248a.de-de.xlf"
$wsD.Range("E3").Value = "2016-03-22 17:12:59"
$wsD.Range("F3").Value = "59db17ff-0d89-493f-a4fb-64bdf414a197.md"
$wsD.Range("G3").Value = "59db17ff-0d89-493f-a4fb-64bdf414a197.bac44bfe36bffe9cc476143af072f2ebdc47248a.de-de.xlf"
$wsD.Range("H3").Value = "2016-03-22 17:13:30"
$wsD.Range("J3").Value = "Include"

$wsD.Range("A4").Value = "198542c6-4d12-4bf5-824b-0062f0c832fd.md"
$wsD.Range("B4").Value = ".md"
$wsD.Range("C4").Value = "Ready for handoff"
$wsD.Range("D4").Value = "198542c6-4d12-4bf5-824b-0062f0c832fd.fa0c9f31959e458e7e6c4b18a2ee6277848d5c38.de-de.xlf"
$wsD.Range("E4").Value = "2016-03-22 17:17:23"
$wsD.Range("F4").Value = "198542c6-4d12-4bf5-824b-0062f0c832fd.md"
$wsD.Range("G4").Value = "198542c6-4d12-4bf5-824b-0062f0c832fd.fa0c9f31959e458e7e6c4b18a2ee6277848d5c38.de-de.xlf"
$wsD.Range("H4").Value = "2016-03-22 17:16:27"
$wsD.Range("J4").Value = "Include"

$wsD.Hyperlinks.Delete()
$wsD.Hyperlinks.Add($wsD.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/85732830159f30e0569f744e3032dbff5009f940/e2e/ffffc5d5a3c0-02f9-4907-ab3a-fb673b2c366d.md", "", "", "ffffc5d5a3c0-02f9-4907-ab3a-fb673b2c366d.md") | Out-Null
$wsD.Hyperlinks.Add($wsD.Range("D2"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/83edba7cfcf10a304a23586e26fbfa94e4a18fcb/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/59db17ff-0d89-493f-a4fb-64bdf414a197.bac44bfe36bffe9cc476143af072f2ebdc47248a.de-de.xlf", "", "", "59db17ff-0d89-493f-a4fb-64bdf414a197.bac44bfe36bffe9cc476143af072f2ebdc47248a.de-de.xlf") | Out-Null
$wsD.Hyperlinks.Add($wsD.Range("F2"), "https://github.com/OpenLocalizationTestOrg/oltest.de-de/blob/d07e56492f8c6d3d6a9d305414f3bbff64e8444f/e2e/59db17ff-0d89-493f-a4fb-64bdf414a197.md", "", "", "59db17ff-0d89-493f-a4fb-64bdf414a197.md") | Out-Null
$wsD.Hyperlinks.Add($wsD.Range("G2"), "https://github.com/OpenLocalizationTestOrg/olhandback/blob/511a4f6c9a951cf83bd11d0316a227963fc310a9/ol-handback/OpenLocalizationTestOrg/oltest.de-de/ci/ht/59db17ff-0d89-493f-a4fb-64bdf414a197.bac44bfe36bffe9cc476143af072f2ebdc47248a.de-de.xlf", "", "", "59db17ff-0d89-493f-a4fb-64bdf414a197.bac44bfe36bffe9cc476143af072f2ebdc47248a.de-de.xlf") | Out-Null

$wsD.Hyperlinks.Add($wsD.Range("A3"), "https://github.com/OpenLocalizationTest/oltest/blob/85732830159f30e0569f744e3032dbff5009f940/e2e/ffffff791f7ac9-e4ac-4cdf-acf3-2ad41c40a99e.md", "", "", "ffffff791f7ac9-e4ac-4cdf-acf3-2ad41c40a99e.md") | Out-Null
$wsD.Hyperlinks.Add($wsD.Range("D3"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/83edba7cfcf10a304a23586e26fbfa94e4a18fcb/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/59db17ff-0d89-493f-a4fb-64bdf414a197.bac44bfe36bffe9cc476143af072f2ebdc47248a.de-de.xlf", "", "", "59db17ff-0d89-493f-a4fb-64bdf414a197.bac44bfe36bffe9cc476143af072f2ebdc47248a.de-de.xlf") | Out-Null
$wsD.Hyperlinks.Add($wsD.Range("F3"), "https://github.com/OpenLocalizationTestOrg/oltest.de-de/blob/d07e56492f8c6d3d6a9d305414f3bbff64e8444f/e2e/59db17ff-0d89-493f-a4fb-64bdf414a197.md", "", "", "59db17ff-0d89-493f-a4fb-64bdf414a197.md") | Out-Null
$wsD.Hyperlinks.Add($wsD.Range("G3"), "https://github.com/OpenLocalizationTestOrg/olhandback/blob/511a4f6c9a951cf83bd11d0316a227963fc310a9/ol-handback/OpenLocalizationTestOrg/oltest.de-de/ci/ht/59db17ff-0d89-493f-a4fb-64bdf414a197.bac44bfe36bffe9cc476143af072f2ebdc47248a.de-de.xlf", "", "", "59db17ff-0d89-493f-a4fb-64bdf414a197.bac44bfe36bffe9cc476143af072f2ebdc47248a.de-de.xlf") | Out-Null

$wsD.Hyperlinks.Add($wsD.Range("A4"), "https://github.com/OpenLocalizationTest/oltest/blob/85732830159f30e0569f744e3032dbff5009f940/e2e/198542c6-4d12-4bf5-824b-0062f0c832fd.md", "", "", "198542c6-4d12-4bf5-824b-0062f0c832fd.md") | Out-Null
$wsD.Hyperlinks.Add($wsD.Range("D4"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/1f39199008aac02bf6d9e20ba5acdbf2fdfc8753/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/198542c6-4d12-4bf5-824b-0062f0c832fd.fa0c9f31959e458e7e6c4b18a2ee6277848d5c38.de-de.xlf", "", "", "198542c6-4d12-4bf5-824b-0062f0c832fd.fa0c9f31959e458e7e6c4b18a2ee6277848d5c38.de-de.xlf") | Out-Null
$wsD.Hyperlinks.Add($wsD.Range("F4"), "https://github.com/OpenLocalizationTestOrg/oltest.de-de/blob/b9d47025bf25f0f8b6c3507ad8c9a09315f256ec/e2e/198542c6-4d12-4bf5-824b-0062f0c832fd.md", "", "", "198542c6-4d12-4bf5-824b-0062f0c832fd.md") | Out-Null
$wsD.Hyperlinks.Add($wsD.Range("G4"), "https://github.com/OpenLocalizationTestOrg/olhandback/blob/dce3b0f6ad2d8e3069fe31f56ecd7ef486651bed/ol-handback/OpenLocalizationTestOrg/oltest.de-de/ci/ht/198542c6-4d12-4bf5-824b-0062f0c832fd.fa0c9f31959e458e7e6c4b18a2ee6277848d5c38.de-de.xlf", "", "", "198542c6-4d12-4bf5-824b-0062f0c832fd.fa0c9f31959e458e7e6c4b18a2ee6277848d5c38.de-de.xlf") | Out-Null
